$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Split the "Figure 1" image caption into two italic runs:
#    "...the rest of the article. " | "They should be placed underneath the figure."
# -----------------------------------------------------------------
$captionText = "Figure 1. A descriptive caption should be given for all figures, understandable without reference to the rest of the article."

$findRange = $d.Content
$found = $findRange.Find.Execute($captionText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # $findRange now spans exactly the matched caption text (Find collapses/
    # resizes the range it was executed on to the hit).
    $paraIndex = $findRange.Paragraphs.Item(1).Index

    # Add the separating space to the end of the existing (first) run.
    $p = $d.Paragraphs.Item($paraIndex)
    $endOfText = $d.Range($p.Range.End - 1, $p.Range.End - 1)
    $endOfText.InsertAfter(" ")

    # Create a new paragraph right after it to hold the second sentence as
    # its own run (a fresh paragraph always starts a fresh run), then fold
    # the paragraph break back out so both runs end up in the same
    # paragraph without Word merging them into a single run.
    $p = $d.Paragraphs.Item($paraIndex)
    $endOfPara = $d.Range($p.Range.End - 1, $p.Range.End - 1)
    $endOfPara.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($paraIndex + 1)
    $newPara.Range.InsertAfter("They should be placed underneath the figure.")

    $p = $d.Paragraphs.Item($paraIndex)
    $markRange = $d.Range($p.Range.End - 1, $p.Range.End)
    $markRange.Delete()
}

# -----------------------------------------------------------------
# 2. Register the new (pandoc reference) paragraph styles, each based
#    on "Normal": Addressee, Horizontal Line, Header Left.
# -----------------------------------------------------------------
$normal = $d.Styles("Normal")

$addressee = $d.Styles.Add("Addressee", 1)
$addressee.BaseStyle = $normal

$horizontalLine = $d.Styles.Add("Horizontal Line", 1)
$horizontalLine.BaseStyle = $normal

$headerLeft = $d.Styles.Add("Header Left", 1)
$headerLeft.BaseStyle = $normal
